$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 570 (shifts existing rows 570-611 down to 571-612)
$ws.Rows.Item(570).Insert()

# Fill in the new row's data: 2026/01/07 (Wed), hour 3, ranking 185
# Force column A to be treated as literal text (matching the existing
# date-as-text cells) rather than being auto-converted to a date serial,
# then clear the temporary number format so no stray style is left behind.
$ws.Range("A570").NumberFormat = "@"
$ws.Range("A570").Value = "2026/01/07"
$ws.Range("A570").ClearFormats()

$ws.Range("B570").Value = "水"
$ws.Range("C570").Value = 3
$ws.Range("D570").Value = 185
